$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New row 31: next task added to the "Solved" tracker table ----

# A31: sequence number (continues numbering from row 30)
$ws.Range("A31").Value = 7

# B31: task description, re-using the wrap-text style already used by the
# other multi-line rows in the table (B20 / B30), then set the new text
$ws.Range("B30").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("B31").Value = "Add Data for testing (50 different users with lots of `nconversations, products, etc)"

# C31: status cell - start from the existing red "Tehnical" style cell
# (red fill + medium border), then additionally color the (empty) font
# red, producing a new combined style (new font + existing fill/border)
$ws.Range("C18").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("C31").Font.Color = 255

$excel.CutCopyMode = 0

# Match the row height used by the other wrapped row (30)
$ws.Rows.Item(31).RowHeight = 29.4

# Update the selection to match the recorded cursor position after editing
$ws.Range("D2").Select()
